{"js": "// The edit: delete the paragraph that reads\n// \"\u73af\u5883\u914d\u7f6e\uff1aGPU\uff1aRTX GeForce 4070super 12GB  CUDA\uff1a12.1  Pytorch\uff1a2.4.1\"\n// (its run + its paragraph mark), leaving the preceding \"\u6ce8\uff1a\" paragraph\n// directly followed by the \"\u6bcf\u4e2a\u6a21\u578bepoch=10\uff0cbatch_size=128...\" paragraph.\n\nconst targetText = \"\u73af\u5883\u914d\u7f6e\uff1aGPU\uff1aRTX GeForce 4070super 12GB  CUDA\uff1a12.1  Pytorch\uff1a2.4.1\";\n\nconst body = context.document.body;\nconst results = body.search(targetText, { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Target paragraph text not found: \" + targetText);\n}\n\n// Resolve the full paragraph that contains the matched text and delete it\n// entirely (this also removes its paragraph mark, merging the surrounding\n// paragraphs' flow without leaving an empty paragraph behind).\nconst paras = results.items[0].paragraphs;\nparas.load(\"items\");\nawait context.sync();\n\nparas.items[0].delete();\nawait context.sync();\n", "ps1": "# The edit: delete the paragraph that reads\n# \"\u73af\u5883\u914d\u7f6e\uff1aGPU\uff1aRTX GeForce 4070super 12GB  CUDA\uff1a12.1  Pytorch\uff1a2.4.1\"\n# (its run + its paragraph mark), leaving the preceding \"\u6ce8\uff1a\" paragraph\n# directly followed by the \"\u6bcf\u4e2a\u6a21\u578bepoch=10\uff0cbatch_size=128...\" paragraph.\n\n$d = $word.ActiveDocument\n\n$targetText = \"\u73af\u5883\u914d\u7f6e\uff1aGPU\uff1aRTX GeForce 4070super 12GB  CUDA\uff1a12.1  Pytorch\uff1a2.4.1\"\n\n$rng = $d.Content\n$found = $rng.Find.Execute($targetText)\n\nif (-not $found) {\n    throw \"Target paragraph text not found: $targetText\"\n}\n\n# Expand the found range to the whole paragraph (this grabs the trailing\n# paragraph mark too), then delete it outright so the paragraph disappears\n# instead of leaving an empty paragraph behind.\n$para = $rng.Paragraphs(1)\n$para.Range.Delete()\n"}
